$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The four "team" rows that used to live in A11:B14 (a timestamp + a team-name
# answer) are being moved up into A2:B5 as the new "next" answers to display,
# and the old slots (A11:B14) are being cleared out, ready for the next
# submissions to come in.

# 1) Copy the existing cell formatting (styles s="2"/s="3") from A11:B14 up to
#    A2:B5 so the moved-in cells keep the same look (date format column / right
#    aligned wrapped text column) without touching A11:B14 itself.
$ws.Range("A11:B14").Copy()
$ws.Range("A2:B5").PasteSpecial(-4122)

# 2) Put the new team-name text into the B column of the newly formatted rows.
$ws.Range("B2").Value = "הקבוצה של: נועם, הקשבי, המפקד"
$ws.Range("B3").Value = "הקבוצה של: טון, אריה"
$ws.Range("B4").Value = "הקבוצה של: אלכס, היייי"
$ws.Range("B5").Value = "הקבוצה של: אור, שרי"

# 3) Row 2 grew a little taller to fit the longer, 3-name answer.
$ws.Rows(2).RowHeight = 27

# 4) Clear out the old answers from A11:B14 (formatting/styles stay as-is).
$ws.Range("A11:B14").ClearContents()

# 5) Move the active selection to B2, where the newest answer now lives.
$ws.Range("B2").Select() | Out-Null
